# Generate Report for Archive
# Updates the "Status" text from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn and de-de sheets, and re-sizes the
# corresponding Status columns (AutoFit side effect of the shorter text).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $text = [string]$cell.Text
        if ($text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# Narrower "Status" columns now that the shorter label fits (matches the
# post-edit autofit width of 13.4101845877511 "characters" as closely as
# this runtime's ColumnWidth setter allows).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
